$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.064.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.021.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.78"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0797"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.323.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.827"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.35"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.023.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.101.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.21"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.26"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.53"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.15"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.82"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.74"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0664"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.51"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.52"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.09%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.33"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0215"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.36"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.398.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.77"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.61%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.06"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +9.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.213.62"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.59%  "
